# Update the pl_mw.xlsx results sheet ("case with 380 kV done"):
# recomputed values for columns B, C, D, F, G, H, I, L, O across rows 2-25.
$data = @{
    2 = @{ "B"=1.271646772049905; "C"=0.3476061929558512; "D"=0.02441764782842881; "F"=0.4217104846387301; "G"=0.2684733239362203; "H"=0.4447382381477354; "I"=0.4137322083167376; "L"=0.2999281063084425; "O"=1.347907406768925 }
    3 = @{ "B"=1.125350326766807; "C"=0.3334892202241519; "D"=0.02159442738581419; "F"=0.4221922092944936; "G"=0.2700995352384936; "H"=0.4495858651692117; "I"=0.4228202327628203; "L"=0.2883149905328111; "O"=1.361238472332232 }
    4 = @{ "B"=1.035251821604675; "C"=0.3248484222071397; "D"=0.01985219986630682; "F"=0.4228811660229113; "G"=0.271418766893099; "H"=0.45284773932147; "I"=0.4287779301504386; "L"=0.281354871719131; "O"=1.370692078904241 }
    5 = @{ "B"=0.9984702208467411; "C"=0.3213344333515806; "D"=0.01914007419814112; "F"=0.4232606183171157; "G"=0.2720367497487999; "H"=0.4542487088830711; "I"=0.4313005042171589; "L"=0.2785614762152306; "O"=1.374862767773067 }
    6 = @{ "B"=0.9923587723892524; "C"=0.320751385176294; "D"=0.01902169762374228; "F"=0.4233295821069873; "G"=0.2721442125982705; "H"=0.4544856699202242; "I"=0.4317250936482591; "L"=0.2781002285967702; "O"=1.37557451088523 }
    7 = @{ "B"=1.034756034688883; "C"=0.3248010016177432; "D"=0.01984260454625542; "F"=0.4228858840566545; "G"=0.2714267761020395; "H"=0.4528663429104327; "I"=0.4288115670231747; "L"=0.2813170251832275; "O"=1.37074703848748 }
    8 = @{ "B"=1.221261724423528; "C"=0.3427332726151064; "D"=0.02344604759914404; "F"=0.4217948761496189; "G"=0.2689673141132545; "H"=0.4463504390574755; "I"=0.4167872724406045; "L"=0.2958885763576831; "O"=1.352240335153724 }
    9 = @{ "B"=1.584743363193184; "C"=0.3780977549428997; "D"=0.03044108195143735; "F"=0.4227834363240106; "G"=0.2667014069140379; "H"=0.4358392208957369; "I"=0.3962120373597653; "L"=0.3258142095225196; "O"=1.326042017033714 }
    10 = @{ "B"=1.850313008249771; "C"=0.4041835376434619; "D"=0.03553491321693514; "F"=0.4254291682932703; "G"=0.2666128248144659; "H"=0.4295009284947469; "I"=0.3829377016223123; "L"=0.3486254481999396; "O"=1.312989789252541 }
    11 = @{ "B"=1.970785659475496; "C"=0.4160695147502338; "D"=0.03784198416842344; "F"=0.4270522851275871; "G"=0.2669184914303457; "H"=0.4269186304793848; "I"=0.377301171570215; "L"=0.3591824440480877; "O"=1.308406466292041 }
    12 = @{ "B"=2.016354895944289; "C"=0.4205728652550818; "D"=0.03871411129863134; "F"=0.4277274458516942; "G"=0.2670842768064006; "H"=0.4259841159530851; "I"=0.3752247758738037; "L"=0.3632059707822464; "O"=1.306866331076918 }
    13 = @{ "B"=2.00654306411775; "C"=0.4196028882912799; "D"=0.03852635113558733; "F"=0.4275793432061761; "G"=0.2670463420475855; "H"=0.4261834516323262; "I"=0.3756693801154984; "L"=0.362338284997719; "O"=1.307189321509924 }
    14 = @{ "B"=1.974535708793326; "C"=0.4164399625283579; "D"=0.03791376516784339; "F"=0.4271066169877287; "G"=0.2669311262605447; "H"=0.4268408783587034; "I"=0.3771291808395034; "L"=0.3595129449098522; "O"=1.308275836735646 }
    15 = @{ "B"=1.954923537131947; "C"=0.4145028783561315; "D"=0.0375383402931817; "F"=0.4268249458064162; "G"=0.2668670779131759; "H"=0.4272492178409308; "I"=0.3780309155608474; "L"=0.3577857031896912; "O"=1.308966836386944 }
    16 = @{ "B"=1.84243283148561; "C"=0.4034071172962399; "D"=0.03538393218543945; "F"=0.4253315526582568; "G"=0.2665998350317409; "H"=0.4296757501456554; "I"=0.383314171213053; "L"=0.3479391408219641; "O"=1.313316633337791 }
    17 = @{ "B"=1.773335303783085; "C"=0.3966049132593525; "D"=0.03405963834643444; "F"=0.4245229992771158; "G"=0.2665247045819186; "H"=0.4312414904349495; "I"=0.3866584072546804; "L"=0.3419446593830173; "O"=1.316332462771882 }
    18 = @{ "B"=1.733560742861528; "C"=0.3926943175536337; "D"=0.03329698841409368; "F"=0.4240974205769135; "G"=0.266514038209607; "H"=0.4321703930517202; "I"=0.388619751167294; "L"=0.3385137434325856; "O"=1.318194508210453 }
    19 = @{ "B"=1.720088436523724; "C"=0.391370587532208; "D"=0.03303860649900514; "F"=0.4239601018279089; "G"=0.2665160078986588; "H"=0.4324897674550314; "I"=0.3892903186298486; "L"=0.3373550077381395; "O"=1.318846826765977 }
    20 = @{ "B"=1.78069413293224; "C"=0.3973288307413725; "D"=0.03420071050850026; "F"=0.4246049837719283; "G"=0.2665293318013724; "H"=0.431071882207668; "I"=0.3862984909568699; "L"=0.3425810282798949; "O"=1.315998229715376 }
    21 = @{ "B"=1.983938454948031; "C"=0.4173689288610944; "D"=0.03809373793055215; "F"=0.4272438241014527; "G"=0.2669636076904283; "H"=0.4266465993363013; "I"=0.3766988251205827; "L"=0.3603421155461461; "O"=1.307951390039335 }
    22 = @{ "B"=2.116470999645287; "C"=0.4304800313677504; "D"=0.04062922734858887; "F"=0.4293213111416634; "G"=0.2675392120402336; "H"=0.4240071000072021; "I"=0.3707632696381697; "L"=0.3721004880366507; "O"=1.303832000016399 }
    23 = @{ "B"=2.045764191462297; "C"=0.4234812644764929; "D"=0.03927681413611595; "F"=0.4281801692753291; "G"=0.2672052090750299; "H"=0.425392710997869; "I"=0.3739001515678986; "L"=0.3658110778247021; "O"=1.30592607275662 }
    24 = @{ "B"=1.777367360163908; "C"=0.3970015473613273; "D"=0.0341369358164485; "F"=0.4245677962754542; "G"=0.2665271385380237; "H"=0.4311484725752379; "I"=0.3864610885899928; "L"=0.3422932779563297; "O"=1.316148937186796 }
    25 = @{ "B"=1.486664389307521; "C"=0.3685112637388102; "D"=0.02855657613197593; "F"=0.4221798263982137; "G"=0.2670388432517186; "H"=0.4384398733895054; "I"=0.4014554316263048; "L"=0.31757377310295; "O"=1.332044382984904 }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
